# Update the "Förändrad" (changed) date column (C) for rows 2-27 from
# 45291 (2023-12-31) to 45292 (2024-01-01), keeping existing cell style/format.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

for ($row = 2; $row -le 27; $row++) {
    $ws.Cells.Item($row, 3).Value = 45292
}
